# Apply updated crypto price/volume figures to columns D (Price) and E (Volume(1h)).
# Some new Price values (e.g. "0.995") look like plain decimal numbers; assigning them
# directly to .Value would make Excel auto-convert the cell to a Number (and e.g. drop the
# trailing zero in "88.50"). Prefixing with a single quote keeps them as text, matching
# the source data which stores every Price/Volume cell as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.810.56'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '1.626.09'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('D4').Value = '''0.995'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = '''210.87'
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('E6').Value = '  -1.00%  '
$ws.Range('D7').Value = '''0.996'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').Value = '''23.15'
$ws.Range('E8').Value = '  -1.63%  '
$ws.Range('E9').Value = '  -0.83%  '
$ws.Range('E10').Value = '  -1.18%  '
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').Value = '1.857.22'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '1.632.27'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('D15').Value = '''0.555'
$ws.Range('E15').Value = '  -1.39%  '
$ws.Range('D16').Value = '''64.87'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').Value = '27.828.57'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '''228.06'
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('E20').Value = '  -1.26%  '
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('E23').Value = '  -4.94%  '
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('D25').Value = '''155.28'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').Value = '''0.111'
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('D28').Value = '''15.44'
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').Value = '1.408.27'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  +2.11%  '
$ws.Range('D36').Value = '''0.997'
$ws.Range('E36').Value = '  -0.66%  '
$ws.Range('E37').Value = '  -1.53%  '
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('E39').Value = '  -1.12%  '
$ws.Range('E40').Value = '  -2.58%  '
$ws.Range('D41').Value = '''0.996'
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('E42').Value = '  -2.20%  '
$ws.Range('D43').Value = '''65.68'
$ws.Range('E43').Value = '  -2.19%  '
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('D45').Value = '''5.42'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('D46').Value = '1.766.31'
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range('E47').Value = '  -3.74%  '
$ws.Range('D48').Value = '''88.50'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('E49').Value = '  +0.94%  '
$ws.Range('E50').Value = '  -0.61%  '
$ws.Range('D51').Value = '''7.57'
$ws.Range('E51').Value = '  +0.25%  '
